$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $text) {
    $stage = $ws.Range("ZZ1")
    $stage.NumberFormat = "@"
    $stage.Value = $text
    $stage.Copy()
    $ws.Range($rangeAddr).PasteSpecial(-4163)
}

$ws.Range("D2").Value = "64.852.45"
$ws.Range("E2").Value = "  +1.35%  "
$ws.Range("D3").Value = "3.172.91"
$ws.Range("E3").Value = "  +1.33%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("E5").Value = "  +1.82%  "
$ws.Range("E6").Value = "  -1.75%  "
$ws.Range("E7").Value = "  +0.20%  "
$ws.Range("D8").Value = "3.158.73"
$ws.Range("E8").Value = "  +0.91%  "
$ws.Range("E9").Value = "  -1.05%  "
$ws.Range("E10").Value = "  -0.23%  "
$ws.Range("E11").Value = "  -1.72%  "
$ws.Range("E12").Value = "  -0.83%  "
$ws.Range("E13").Value = "  -0.29%  "
$ws.Range("E14").Value = "  -2.99%  "
$ws.Range("D15").Value = "3.693.19"
$ws.Range("E15").Value = "  +2.11%  "
$ws.Range("E16").Value = "  +2.76%  "
$ws.Range("D17").Value = "64.824.58"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "3.171.95"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("E19").Value = "  -1.37%  "
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -0.02%  "
$ws.Range("E27").Value = "  -3.08%  "
$ws.Range("E28").Value = "  -0.30%  "
$ws.Range("E29").Value = "  -6.49%  "
$ws.Range("E30").Value = "  -2.20%  "
$ws.Range("E31").Value = "  -7.88%  "
$ws.Range("E32").Value = "  -0.02%  "
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("E34").Value = "  -0.51%  "
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("D36").Value = "0.0₃0783"
$ws.Range("E36").Value = "  +3.68%  "
$ws.Range("E37").Value = "  -1.56%  "
$ws.Range("E38").Value = "  -1.22%  "
$ws.Range("E39").Value = "  -2.47%  "
$ws.Range("E40").Value = "  +2.91%  "
$ws.Range("E41").Value = "  +0.00%  "
$ws.Range("E42").Value = "  -3.43%  "
$ws.Range("E43").Value = "  -1.47%  "
$ws.Range("D44").Value = "2.853.96"
$ws.Range("E44").Value = "  -1.14%  "
$ws.Range("E45").Value = "  -2.13%  "
$ws.Range("E46").Value = "  -2.05%  "
$ws.Range("E47").Value = "  +4.79%  "
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("E49").Value = "  +0.13%  "
$ws.Range("E50").Value = "  -1.58%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("E51").Value = "  +0.90%  "

Set-TextValue "D5" "615.66"
Set-TextValue "D6" "147.17"
Set-TextValue "D10" "0.153"
Set-TextValue "D11" "5.51"
Set-TextValue "D12" "0.475"
Set-TextValue "D20" "478.79"
Set-TextValue "D21" "14.71"
Set-TextValue "D24" "13.81"
Set-TextValue "D25" "84.66"
Set-TextValue "D27" "2.84"
Set-TextValue "D28" "8.63"
Set-TextValue "D30" "6.91"
Set-TextValue "D31" "2.08"
Set-TextValue "D34" "26.64"
Set-TextValue "D38" "3.23"
Set-TextValue "D39" "53.17"
Set-TextValue "D40" "465.40"
Set-TextValue "D41" "0.0401"
Set-TextValue "D45" "2.34"
Set-TextValue "D46" "0.269"
Set-TextValue "D47" "2.44"
Set-TextValue "D48" "26.70"
Set-TextValue "D51" "120.76"

$ws.Range("ZZ1").Clear()
